# Applies the "Add files via upload" edit to the Team_Management sheet:
#  - Column D (Rig_No.) values are changed from the generic A/B/C placeholders
#    to the actual rig names (or "Not Determined" where unknown).
#  - Columns E (Job Type), F (Actual Date) and G (Spent Days) are cleared
#    for every data row, since that information is no longer tracked here.
#  - Row 8's Team_No. is corrected from Team_2 to Team_3.
#  - The worksheet's active selection moves to D20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team_Management")

# Fix the team number on row 8 (was Team_2, should be Team_3)
$ws.Range("A8").Value = "Team_3"

# New values for column D (Rig_No.) for rows 2-19
$rigNo = @{
    2  = "EDC-12"
    3  = "Home"
    4  = "Not Determined"
    5  = "EDC-92"
    6  = "Home"
    7  = "Not Determined"
    8  = "EDC-40"
    9  = "Home"
    10 = "Not Determined"
    11 = "Home"
    12 = "Not Determined"
    13 = "EDC-54"
    14 = "EDC-82"
    15 = "EDC-88"
    16 = "EDC-12"
    17 = "Home"
    18 = "Not Determined"
    19 = "HT-101"
}

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 4).Value = $rigNo[$row]
    $ws.Cells.Item($row, 5).ClearContents()
    $ws.Cells.Item($row, 6).ClearContents()
    $ws.Cells.Item($row, 7).ClearContents()
}

# Move the active selection to D20, matching the saved view state
$ws.Range("D20").Select()

$wb.Save()
